$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Isi Dunia")

# DNA column updates (col D) for Polisi, Pemburu, Rumput, Pohon rows
$ws.Range("D5").Value = "C"
$ws.Range("D6").Value = "H"
$ws.Range("D11").Value = "^"
$ws.Range("D12").Value = "!"

# Rename constructor-call demo label from Burung_Unta to Hewan
$ws.Range("C16").Value = "Hewan"

# Add missing first-letter cells for the Gajah / Burung_Unta examples
$ws.Range("E21").Value = "b"
$ws.Range("E23").Value = "b"

# Scroll / selection state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("E19").Select()
